# ---------------------------------------------------------------------------
# testData.xlsx edit script
#
# Summary of changes (per commit "little changes done in framework constants
# files 18-10"):
#   * RETAIL_DATA becomes the active/selected sheet (activeTab=1).
#   * RUNMANAGER (sheet 1) header row renamed to nicer display text and the
#     "execute" flag for the first data row flips from yes -> no. Selection
#     anchor moves from C10 to B10.
#   * RETAIL_DATA (sheet 2) gets a new sub-header row describing the last two
#     (data) columns, a renamed header row, an extra (unused) column F, and
#     the trailing blank styled rows are trimmed away. Selection anchor moves
#     to F4.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsRun    = $wb.Worksheets.Item(1)   # RUNMANAGER
$wsRetail = $wb.Worksheets.Item(2)   # RETAIL_DATA

# ---------------------------------------------------------------------------
# 1. RUNMANAGER sheet
# ---------------------------------------------------------------------------

# Header row - nicer display labels.
$wsRun.Range("A1").Value2 = "Test Case Name"
$wsRun.Range("B1").Value2 = "Test Description"
$wsRun.Range("C1").Value2 = "Execute"
$wsRun.Range("D1").Value2 = "Priority"
$wsRun.Range("E1").Value2 = "Count"

# First data row: execute flag flips to "no".
$wsRun.Range("C2").Value2 = "no"

# Move the saved selection anchor.
[void]$wsRun.Range("B10").Select()

# ---------------------------------------------------------------------------
# 2. RETAIL_DATA sheet
# ---------------------------------------------------------------------------

# Insert a new row 2 (existing rows 2 & 3 shift down to 3 & 4), then drop the
# old trailing placeholder rows (formerly 4-7, now shifted to 5-8).
[void]$wsRetail.Rows.Item(2).Insert()
[void]$wsRetail.Range("A5:E8").EntireRow.Delete()

# Header row - nicer display labels, plus two new "data" columns replacing
# the old username/password headers.
$wsRetail.Range("A1").Value2 = "Test Case Name"
$wsRetail.Range("B1").Value2 = "Execute"
$wsRetail.Range("C1").Value2 = "Browser"
$wsRetail.Range("D1").Value2 = "data1"
$wsRetail.Range("E1").Value2 = "data2"

# New sub-header row: A2:C2 stay visually blank (quote-prefixed empty text,
# matching the style already used for the "1" cells on RUNMANAGER), D2/E2
# describe what data1/data2 actually hold.
$wsRetail.Range("A2").Formula = "'"
$wsRetail.Range("B2").Formula = "'"
$wsRetail.Range("C2").Formula = "'"
$wsRetail.Range("D2").Value2 = "User name"
$wsRetail.Range("E2").Value2 = "Login Password"

# Re-point the existing hyperlinks at the cells they now land on (E3/E4) -
# the underlying target address is unchanged.
[void]$wsRetail.Hyperlinks.Delete()
[void]$wsRetail.Hyperlinks.Add($wsRetail.Range("E3"), "mailto:Asdf@123")
[void]$wsRetail.Hyperlinks.Add($wsRetail.Range("E4"), "mailto:Asdf@123")

# New (currently unused) column F, sized to fit a "Login Password"-ish label.
$wsRetail.Columns.Item(6).ColumnWidth = 15.14

# Move the saved selection anchor and make RETAIL_DATA the active sheet/tab.
[void]$wsRetail.Range("F4").Select()
[void]$wsRetail.Activate()
